$updates = @(
    @{ Row = 2; D = "28.038.77"; E = "  -0.15%  " }
    @{ Row = 3; D = "1.871.77"; E = "  -1.07%  " }
    @{ Row = 4; D = "1.005"; E = "  +0.40%  " }
    @{ Row = 5; D = "312.75"; E = "  -0.38%  " }
    @{ Row = 6; D = "1.004"; E = "  +0.32%  " }
    @{ Row = 7; D = "0.5086"; E = "  +0.69%  " }
    @{ Row = 8; D = "0.3816"; E = "  -2.06%  " }
    @{ Row = 9; D = $null; E = "  -10.52%  " }
    @{ Row = 10; D = "1.113"; E = "  -1.35%  " }
    @{ Row = 11; D = "41.55"; E = "  -0.61%  " }
    @{ Row = 12; D = "6.224"; E = "  -2.53%  " }
    @{ Row = 13; D = "1.870.81"; E = "  -1.31%  " }
    @{ Row = 14; D = "20.50"; E = "  -1.49%  " }
    @{ Row = 15; D = "7.206"; E = "  -1.15%  " }
    @{ Row = 16; D = "1.005"; E = "  +0.40%  " }
    @{ Row = 17; D = $null; E = "  -1.06%  " }
    @{ Row = 18; D = "90.90"; E = "  -1.45%  " }
    @{ Row = 19; D = "0.06636"; E = "  -0.25%  " }
    @{ Row = 20; D = "17.97"; E = "  +0.67%  " }
    @{ Row = 21; D = $null; E = "  +0.22%  " }
    @{ Row = 22; D = "6.045"; E = "  -2.66%  " }
    @{ Row = 23; D = "28.078.32"; E = "  -0.20%  " }
    @{ Row = 24; D = "11.15"; E = "  -2.18%  " }
    @{ Row = 25; D = "2.266"; E = "  -2.36%  " }
    @{ Row = 26; D = "2.593"; E = "  +1.99%  " }
    @{ Row = 27; D = "2.095.35"; E = "  -0.99%  " }
    @{ Row = 28; D = "157.16"; E = "  -0.87%  " }
    @{ Row = 29; D = "20.58"; E = "  -1.19%  " }
    @{ Row = 30; D = "125.78"; E = "  -0.91%  " }
    @{ Row = 31; D = "0.1055"; E = "  +0.03%  " }
    @{ Row = 32; D = "1.045"; E = "  -2.98%  " }
    @{ Row = 33; D = "5.603"; E = "  -0.08%  " }
    @{ Row = 34; D = "3.606"; E = "  -0.07%  " }
    @{ Row = 35; D = "9.682"; E = "  +2.30%  " }
    @{ Row = 36; D = $null; E = "  +1.94%  " }
    @{ Row = 37; D = "0.06590"; E = "  -0.25%  " }
    @{ Row = 38; D = "0.2167"; E = "  -1.39%  " }
    @{ Row = 39; D = $null; E = "  -0.31%  " }
    @{ Row = 40; D = "0.6479"; E = "  +0.63%  " }
    @{ Row = 41; D = $null; E = "  -7.18%  " }
    @{ Row = 42; D = "11.33"; E = "  -2.88%  " }
    @{ Row = 43; D = "4.884"; E = "  -1.73%  " }
    @{ Row = 44; D = "0.6127"; E = "  +1.26%  " }
    @{ Row = 45; D = "13.13"; E = "  -1.14%  " }
    @{ Row = 46; D = "1.296"; E = "  -0.33%  " }
    @{ Row = 47; D = "3.666"; E = "  -0.58%  " }
    @{ Row = 48; D = "2.011"; E = "  +0.42%  " }
    @{ Row = 49; D = $null; E = "  +2.01%  " }
    @{ Row = 50; D = "120.85"; E = "  -1.03%  " }
    @{ Row = 51; D = "80.40"; E = "  +1.95%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$($u.Row)")
        $dVal = $u.D
        if ($dVal -match '^[+-]?[0-9]*\.?[0-9]+$') {
            # Value looks like a plain number (e.g. "1.005") - force the
            # cell to stay text by pre-formatting it as Text before the
            # value is assigned (otherwise Excel silently converts it to
            # a numeric value).
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
